# Update the "想去人数" (want-to-go count, column F) values that were
# refreshed by +1 on the source site for a handful of events, across the
# "展览", "演出" and "全部类型" sheets ("本地生活" is untouched).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 85
$ws1.Range("F6").Value  = 1053
$ws1.Range("F12").Value = 1671
$ws1.Range("F21").Value = 606
$ws1.Range("F22").Value = 12252
$ws1.Range("F23").Value = 12296

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 29

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 85
$ws4.Range("F7").Value  = 1053
$ws4.Range("F13").Value = 1671
$ws4.Range("F25").Value = 606
$ws4.Range("F26").Value = 12252
$ws4.Range("F27").Value = 12296
$ws4.Range("F40").Value = 29
